$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4, 2).Value = 2090115
$ws.Cells.Item(4, 3).Value = 414
$ws.Cells.Item(4, 4).Value = 816457
$ws.Cells.Item(4, 5).Value = 1157606
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 18
$ws.Cells.Item(4, 8).Value = 116052

# Row 7: India -> India
$ws.Cells.Item(7, 2).Value = 298482
$ws.Cells.Item(7, 3).Value = 199
$ws.Cells.Item(7, 4).Value = 147544
$ws.Cells.Item(7, 5).Value = 142426
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 11
$ws.Cells.Item(7, 8).Value = 8512

# Row 12: Alemania -> Alemania
$ws.Cells.Item(12, 2).Value = 186867
$ws.Cells.Item(12, 3).Value = 72
$ws.Cells.Item(12, 4).Value = 171600
$ws.Cells.Item(12, 5).Value = 6415
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 8852

# Row 26: Bielorrusia -> Bielorrusia
$ws.Cells.Item(26, 2).Value = 52520
$ws.Cells.Item(26, 3).Value = 704
$ws.Cells.Item(26, 4).Value = 27760
$ws.Cells.Item(26, 5).Value = 24462
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 5
$ws.Cells.Item(26, 8).Value = 298

# Row 27: Suecia -> Paises Bajos
$ws.Cells.Item(27, 1).Value = "Paises Bajos"
$ws.Cells.Item(27, 2).Value = 48461
$ws.Cells.Item(27, 3).Value = 210
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 9
$ws.Cells.Item(27, 8).Value = 6053

# Row 28: Paises Bajos -> Suecia
$ws.Cells.Item(28, 1).Value = "Suecia"
$ws.Cells.Item(28, 2).Value = 48288
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 4814

# Row 56: Kazajistan -> Kazajistan
$ws.Cells.Item(56, 2).Value = 13872
$ws.Cells.Item(56, 3).Value = 314
$ws.Cells.Item(56, 4).Value = 8668
$ws.Cells.Item(56, 5).Value = 5136
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 1
$ws.Cells.Item(56, 8).Value = 68

# Row 60: Moldavia -> Ghana
$ws.Cells.Item(60, 1).Value = "Ghana"
$ws.Cells.Item(60, 2).Value = 10856
$ws.Cells.Item(60, 3).Value = 498
$ws.Cells.Item(60, 4).Value = 3921
$ws.Cells.Item(60, 5).Value = 6887
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 48

# Row 61: Argelia -> Moldavia
$ws.Cells.Item(61, 1).Value = "Moldavia"
$ws.Cells.Item(61, 2).Value = 10727
$ws.Cells.Item(61, 3).Value = 0
$ws.Cells.Item(61, 4).Value = 6229
$ws.Cells.Item(61, 5).Value = 4116
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 7
$ws.Cells.Item(61, 8).Value = 382

# Row 62: Ghana -> Argelia
$ws.Cells.Item(62, 1).Value = "Argelia"
$ws.Cells.Item(62, 2).Value = 10589
$ws.Cells.Item(62, 3).Value = 0
$ws.Cells.Item(62, 4).Value = 7255
$ws.Cells.Item(62, 5).Value = 2593
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 0
$ws.Cells.Item(62, 8).Value = 741

# Row 76: Tayikistan -> Uzbekistan
$ws.Cells.Item(76, 1).Value = "Uzbekistan"
$ws.Cells.Item(76, 2).Value = 4837
$ws.Cells.Item(76, 3).Value = 96
$ws.Cells.Item(76, 4).Value = 3700
$ws.Cells.Item(76, 5).Value = 1118
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 19

# Row 77: Uzbekistan -> Tayikistan
$ws.Cells.Item(77, 1).Value = "Tayikistan"
$ws.Cells.Item(77, 2).Value = 4834
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 3062
$ws.Cells.Item(77, 5).Value = 1723
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 49

# Row 79: Costa de Marfil -> Republica de Yibuti
$ws.Cells.Item(79, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(79, 2).Value = 4441
$ws.Cells.Item(79, 3).Value = 43
$ws.Cells.Item(79, 4).Value = 2730
$ws.Cells.Item(79, 5).Value = 1673
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = 38

# Row 80: Republica de Yibuti -> Costa de Marfil
$ws.Cells.Item(80, 1).Value = "Costa de Marfil"
$ws.Cells.Item(80, 2).Value = 4404
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(80, 4).Value = 2212
$ws.Cells.Item(80, 5).Value = 2151
$ws.Cells.Item(80, 6).Value = 0
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 41

# Row 85: Republica de Macedonia -> Republica de Macedonia
$ws.Cells.Item(85, 2).Value = 3701
$ws.Cells.Item(85, 3).Value = 163
$ws.Cells.Item(85, 4).Value = 1694
$ws.Cells.Item(85, 5).Value = 1836
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 2
$ws.Cells.Item(85, 8).Value = 171

# Row 92: Bosnia y Herzegovina -> Bosnia y Herzegovina
$ws.Cells.Item(92, 2).Value = 2893
$ws.Cells.Item(92, 3).Value = 61
$ws.Cells.Item(92, 4).Value = 2119
$ws.Cells.Item(92, 5).Value = 611
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 2
$ws.Cells.Item(92, 8).Value = 163

# Row 96: Croacia -> Croacia
$ws.Cells.Item(96, 2).Value = 2249
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 2133
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 107

# Row 114: Albania -> Libano
$ws.Cells.Item(114, 1).Value = "Libano"
$ws.Cells.Item(114, 2).Value = 1422
$ws.Cells.Item(114, 3).Value = 20
$ws.Cells.Item(114, 4).Value = 853
$ws.Cells.Item(114, 5).Value = 538
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 31

# Row 115: Libano -> Albania
$ws.Cells.Item(115, 1).Value = "Albania"
$ws.Cells.Item(115, 2).Value = 1416
$ws.Cells.Item(115, 3).Value = 31
$ws.Cells.Item(115, 4).Value = 1034
$ws.Cells.Item(115, 5).Value = 346
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 36

# Row 133: Georgia -> Georgia
$ws.Cells.Item(133, 2).Value = 843
$ws.Cells.Item(133, 3).Value = 12
$ws.Cells.Item(133, 4).Value = 697
$ws.Cells.Item(133, 5).Value = 133
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 13

# Row 136: San Marino -> San Marino
$ws.Cells.Item(136, 2).Value = 694
$ws.Cells.Item(136, 3).Value = 3
$ws.Cells.Item(136, 4).Value = 520
$ws.Cells.Item(136, 5).Value = 132
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 42

# Row 157: Vietnam -> Vietnam
$ws.Cells.Item(157, 2).Value = 332
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 323
$ws.Cells.Item(157, 5).Value = 9
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 0
